$wb = $excel.ActiveWorkbook
$new = $wb.Worksheets.Add()
try {
  $shp = $new.Shapes.AddTextbox(1, 10, 10, 200, 50)
  $shp.TextFrame.Characters().Text = "Hello World"
  Write-Output "AddTextbox ok"
} catch {
  Write-Output ("AddTextbox failed: " + $_.Exception.Message)
}
try {
  $pic = $new.Shapes.AddPicture("image1.jpg", 0, 1, 10, 100, 100, 100)
  Write-Output "AddPicture ok"
} catch {
  Write-Output ("AddPicture failed: " + $_.Exception.Message)
}
